# Generate Report for Archive
# Re-sort the three localization-status rows (01c3da70…, 72b21453…, e2734160…)
# on every sheet so that the e2734160 entry (still "In Translation") now
# sorts ahead of the two files that are already "Ready for handoff".
#
# The row that used to be on row 5 moves to row 6, the row that used to be
# on row 6 moves to row 7, and the row that used to be on row 7 (e2734160)
# moves up to row 5 — picking up the "In Translation" status it actually
# has instead of the stale "Ready for handoff" that had been left in the
# report.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A5").Value = "e2734160-a6dc-4b16-aa15-d40346028814.md"
$wsOverview.Range("B5").Value = "In Translation"
$wsOverview.Range("C5").Value = "In Translation"
$wsOverview.Range("D5").Value = "2016-31-12 18:31:39"

$wsOverview.Range("A6").Value = "01c3da70-f464-4d60-973a-d00a275bd8ed.md"
$wsOverview.Range("B6").Value = "Ready for handoff"
$wsOverview.Range("C6").Value = "Ready for handoff"
$wsOverview.Range("D6").Value = "2016-29-12 18:29:52"

$wsOverview.Range("A7").Value = "72b21453-9460-4dd3-b944-2e553c742a9f.md"
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-32-12 18:32:02"

# ---------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status |
#              Latest Handoff File | Latest Handoff Datetime | ...
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A5").Value = "e2734160-a6dc-4b16-aa15-d40346028814.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("D5").Value = "e2734160-a6dc-4b16-aa15-d40346028814.cd211a198dd778b7e9dc4b0f15cc7eb7145908be.zh-cn.xlf"
$wsZhCn.Range("E5").Value = "2016-03-12 18:31:36"

$wsZhCn.Range("A6").Value = "01c3da70-f464-4d60-973a-d00a275bd8ed.md"
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "01c3da70-f464-4d60-973a-d00a275bd8ed.121570c80681fa13a7a576a4fe2b81c191d3d534.zh-cn.xlf"
$wsZhCn.Range("E6").Value = "2016-03-12 18:29:48"

$wsZhCn.Range("A7").Value = "72b21453-9460-4dd3-b944-2e553c742a9f.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "72b21453-9460-4dd3-b944-2e553c742a9f.2fa01d6e052989878f19d80f992a4abf8ff4ce80.zh-cn.xlf"
$wsZhCn.Range("E7").Value = "2016-03-12 18:31:58"

# ---------------------------------------------------------------------
# de-de sheet: Source File Name | File Extension | Status |
#              Latest Handoff File | Latest Handoff Datetime | ...
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A5").Value = "e2734160-a6dc-4b16-aa15-d40346028814.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("D5").Value = "e2734160-a6dc-4b16-aa15-d40346028814.cd211a198dd778b7e9dc4b0f15cc7eb7145908be.de-de.xlf"
$wsDeDe.Range("E5").Value = "2016-03-12 18:31:39"

$wsDeDe.Range("A6").Value = "01c3da70-f464-4d60-973a-d00a275bd8ed.md"
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "01c3da70-f464-4d60-973a-d00a275bd8ed.121570c80681fa13a7a576a4fe2b81c191d3d534.de-de.xlf"
$wsDeDe.Range("E6").Value = "2016-03-12 18:29:52"

$wsDeDe.Range("A7").Value = "72b21453-9460-4dd3-b944-2e553c742a9f.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "72b21453-9460-4dd3-b944-2e553c742a9f.2fa01d6e052989878f19d80f992a4abf8ff4ce80.de-de.xlf"
$wsDeDe.Range("E7").Value = "2016-03-12 18:32:02"
